# Modify And Export Network Implemented
#
# The TriangleData sheet's header labels are changed from legacy upper-case
# names to the new lower-case "export" naming convention:
#   Node table   (row 1)  : N_ID/N_X/N_Y          -> n_id/n_x/n_y
#   Edge table   (row 15) : E_ID/N_FROM/N_TO       -> e_id/n_from/n_to
#   Meta table   (row 37) : META_NET_PROP/META_NET_PROP_VALUE -> meta_key/meta_value
#
# (The meta-data values below row 37 - network_size/100x100/total_tension/
#  type/static - and everything on "Sheet1" stay exactly as they were.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TriangleData")

# -- Node header --
$ws.Range("A1").Value = "n_id"
$ws.Range("B1").Value = "n_x"
$ws.Range("C1").Value = "n_y"

# -- Edge header --
$ws.Range("A15").Value = "e_id"
$ws.Range("B15").Value = "n_from"
$ws.Range("C15").Value = "n_to"

# -- Meta header --
$ws.Range("A37").Value = "meta_key"
$ws.Range("B37").Value = "meta_value"

# -- Update the view state left behind by the editing session --
$ws.Activate()
$ws.Range("C37").Select()
$excel.ActiveWindow.Zoom = 130
